# "adds special characters. removes some logs."
#
# 1. Rename the first sheet ("Planilha1" -> "letters").
# 2. Update the data grid on that sheet (rows 34-41, cols AM-AW) to reflect
#    the newly-added "special characters" rows/marks, clearing a couple of
#    stray marks ("logs") along the way.
# 3. Move the active-sheet/selection from "bar_chart" (T16) over to the
#    renamed "letters" sheet (AV43), matching the tab the author left
#    selected when they saved.

$wb = $excel.ActiveWorkbook

# --- 1. rename the sheet -----------------------------------------------
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Name = "letters"

# --- 2. cell edits on the renamed sheet ---------------------------------
$ws.Range("AM34").Value = 1

$ws.Range("AN35").Value = ""

$ws.Range("AM36").Value = 1

$ws.Range("AN37").Value = 1
$ws.Range("AO37").Value = ""

$ws.Range("AN38").Value = 1
$ws.Range("AO38").Value = ""

$ws.Range("AM39").Value = 1
$ws.Range("AP39").Value = 1
$ws.Range("AQ39").Value = 1
$ws.Range("AR39").Value = " "
$ws.Range("AS39").Value = 1
$ws.Range("AT39").Value = 1
$ws.Range("AV39").Value = 1
$ws.Range("AW39").Value = 1

$ws.Range("AN40").Value = ""
$ws.Range("AP40").Value = 1
$ws.Range("AQ40").Value = 1
$ws.Range("AR40").Value = " "
$ws.Range("AS40").Value = 1
$ws.Range("AT40").Value = 1
$ws.Range("AV40").Value = 1
$ws.Range("AW40").Value = 1

$ws.Range("AM41").Value = 1

# --- 3. switch the active tab/selection ---------------------------------
$ws.Activate() | Out-Null
$ws.Range("AV43").Select() | Out-Null
